# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.895.90"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "'1.889.22"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'0.7635"
$ws.Range("E5").Value = "  -1.51%  "

# Row 6
$ws.Range("D6").Value = "'242.71"
$ws.Range("E6").Value = "  -0.65%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.3128"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").Value = "'25.67"
$ws.Range("E9").Value = "  +1.39%  "

# Row 10
$ws.Range("D10").Value = "'0.07158"
$ws.Range("E10").Value = "  -3.34%  "

# Row 11
$ws.Range("D11").Value = "'0.08512"
$ws.Range("E11").Value = "  +4.47%  "

# Row 12
$ws.Range("D12").Value = "'0.7623"
$ws.Range("E12").Value = "  -0.43%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.898.26"
$ws.Range("E13").Value = "  +4.01%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.371"
$ws.Range("E14").Value = "  -1.98%  "

# Row 15
$ws.Range("D15").Value = "'93.84"
$ws.Range("E15").Value = "  +1.43%  "

# Row 16
$ws.Range("D16").Value = "'6.131"
$ws.Range("E16").Value = "  -1.48%  "

# Row 17
$ws.Range("D17").Value = "'29.784.65"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("E18").Value = "  -1.31%  "

# Row 19
$ws.Range("D19").Value = "'243.87"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20
$ws.Range("E20").Value = "  -0.94%  "

# Row 21
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("D22").Value = "'8.006"
$ws.Range("E22").Value = "  -1.39%  "

# Row 23
$ws.Range("D23").Value = "'2.098.75"
$ws.Range("E23").Value = "  -1.02%  "

# Row 24
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "'0.1617"
$ws.Range("E25").Value = "  +2.98%  "

# Row 26
$ws.Range("E26").Value = "  -0.27%  "

# Row 27
$ws.Range("D27").Value = "'161.95"
$ws.Range("E27").Value = "  -0.18%  "

# Row 28
$ws.Range("D28").Value = "'18.78"
$ws.Range("E28").Value = "  -0.15%  "

# Row 29
$ws.Range("E29").Value = "  -0.29%  "

# Row 30
$ws.Range("D30").Value = "'1.484"
$ws.Range("E30").Value = "  +1.97%  "

# Row 31
$ws.Range("D31").Value = "'1.530"
$ws.Range("E31").Value = "  -0.97%  "

# Row 32
$ws.Range("D32").Value = "'4.486"
$ws.Range("E32").Value = "  -0.41%  "

# Row 33
$ws.Range("D33").Value = "'4.093"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34
$ws.Range("D34").Value = "'0.05447"
$ws.Range("E34").Value = "  -2.83%  "

# Row 35
$ws.Range("D35").Value = "'1.244"
$ws.Range("E35").Value = "  -0.58%  "

# Row 36
$ws.Range("D36").Value = "'0.7426"
$ws.Range("E36").Value = "  -2.08%  "

# Row 37
$ws.Range("D37").Value = "'0.9990"

# Row 38
$ws.Range("E38").Value = "  +1.85%  "

# Row 39
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("D40").Value = "'2.779"
$ws.Range("E40").Value = "  -0.44%  "

# Row 41
$ws.Range("D41").Value = "'0.4464"
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$ws.Range("D42").Value = "'1.099.99"
$ws.Range("E42").Value = "  -4.12%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'72.98"
$ws.Range("E43").Value = "  -1.79%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.067"
$ws.Range("E44").Value = "  +1.55%  "

# Row 45
$ws.Range("D45").Value = "'0.8505"
$ws.Range("E45").Value = "  -0.40%  "

# Row 46
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").Value = "'102.88"
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("D48").Value = "'1.869"
$ws.Range("E48").Value = "  -2.02%  "

# Row 49
$ws.Range("D49").Value = "'7.634"
$ws.Range("E49").Value = "  +1.51%  "

# Row 50
$ws.Range("D50").Value = "'2.997"
$ws.Range("E50").Value = "  -4.43%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "'2.003.72"
$ws.Range("E51").Value = "  -1.04%  "

# Clear the quote-prefix formatting residue left by the text-forcing apostrophes
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
